# #5: property aircraft done
# Fix mislabeled property_category values:
#  - 建物 (building) sheet: rows were tagged "land" instead of "building"
#  - 汽車 (car) sheet: row was tagged "land" instead of "car"

$wb = $excel.ActiveWorkbook

$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"

$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
